$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("server")

# Remove the "认证方式" (auth type) column (old column P / index 16).
# This shifts 用户密码/SSH端口/宿主机 left by one column and fixes up
# the shared-string indices, data validations and dimension automatically.
$ws.Columns.Item(16).Delete()

# Add a new data row (row 5) describing asset JR002, mirroring row 2's
# "server" record (JR001) but with its own asset id / management IP.
$ws.Range("A5").Value2 = "server"
$ws.Range("B5").Value2 = "JR002"
$ws.Range("D5").Value2 = 1
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = "10.1.19.11"
$ws.Range("G5").Value2 = 2
$ws.Range("H5").Value2 = 1
$ws.Range("I5").Value2 = 1
$ws.Range("J5").Value2 = 39668
$ws.Range("K5").Value2 = 43320
$ws.Range("L5").Value2 = 800
$ws.Range("N5").Value2 = 0
$ws.Range("O5").Value2 = "root"
$ws.Range("P5").Value2 = "123456"
$ws.Range("Q5").Value2 = 22

# Match the date-style / text-style formatting used by the other rows
# (copy formats only, so no new number-format styles get introduced).
$ws.Range("J2").Copy()
$ws.Range("J5").PasteSpecial(-4122)
$ws.Range("K2").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("P2").Copy()
$ws.Range("P5").PasteSpecial(-4122)

# Fix up the selected cell / scroll position recorded in the sheet view.
$ws.Range("L12").Select()
